$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BatchSheet")

# Update shared strings that currently read "P2P.*" to "Cloud.*"
$ws.Range("E2").Value = "Cloud.createItemBasedRequisition"
$ws.Range("E3").Value = "Cloud.createPurchaseOrder"
$ws.Range("E8").Value = "Cloud.createReceivingReceipt"
$ws.Range("E9").Value = "Cloud.createPurOrderMatchedInvoice"
$ws.Range("E10").Value = "Cloud.createPaymentQuickCheck"

# Update the active sheet's selection from E17 to C16
$ws.Activate()
$ws.Range("C16").Select()
